$d = $word.ActiveDocument

# Locate the "Requisitos" paragraph that lists LOT2056 as a requirement.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOT2056: Trabalho de Conclusão de Curso I (Requisito)*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the LOT2056 requirement paragraph"
}

# Immediately after it the page used to have three trailing paragraphs that
# came from the scraped site chrome:
#   1) an empty spacer paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 ... Jekyll and Github pages ..." footer paragraph
# Those three paragraphs should be removed, leaving the LOT2056 paragraph
# directly followed by the (remaining) trailing empty paragraph and the
# page-break paragraph.
$p1 = $target.Next()
$p2 = $p1.Next()
$p3 = $p2.Next()

$start = $p1.Range.Start
$end = $p3.Range.End
$r = $d.Range($start, $end)
$r.Delete()
